# Add "Implementation" section (rows 21-38) to the Task Assignments -
# Requirements sheet, mirroring the existing REQUIREMENTS / DESIGN
# sections below them (blank separator row, bold section header, then
# alternating-style Task/Assigned To/Completed By rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Seed cell formatting (borders/fill/font) for the new rows by
#    copying it from the existing analogous rows. Rows 11-20 already
#    carry the blank-separator / bold-header / alternating 3-1 style
#    pattern we need for rows 21-30; rows 13-20 carry the plain
#    alternating 3-1 pattern needed for rows 31-38.
# ---------------------------------------------------------------------
$ws.Range("A11:C20").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null

$ws.Range("A13:C20").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null

# A few cells in the copied-from-life data carried an explicit black
# font color (distinct from the rest of the sheet's theme-color font).
# Reproduce that by overriding the font color on just those cells.
$ws.Range("B25:C25").Font.Color = 0
$ws.Range("B26:C26").Font.Color = 0
$ws.Range("C27").Font.Color = 0

# ---------------------------------------------------------------------
# 2) Fill in the cell values, left-to-right / top-to-bottom so shared
#    strings are introduced in the same order the source workbook used.
# ---------------------------------------------------------------------
$ws.Range("A22").Value = "Implementation"

$ws.Range("A23").Value = "Updated Task Requirements and Project Requirements"
$ws.Range("B23").Value = "Will Maberry"
$ws.Range("C23").Value = "Will Maberry"

$ws.Range("A24").Value = "WordList.txt"
$ws.Range("B24").Value = "Will Maberry"
$ws.Range("C24").Value = "Will Maberry"

$ws.Range("A25").Value = "WordLocation.java (and associated unit tests)"
$ws.Range("B25").Value = "Will Maberry"
$ws.Range("C25").Value = "Will Maberry"

$ws.Range("A26").Value = "WordBank.java (and associated unit tests)"
$ws.Range("B26").Value = "Will Maberry"
$ws.Range("C26").Value = "Will Maberry"

$ws.Range("A27").Value = "GameSession.java (and associated unit tests)"
$ws.Range("B27").Value = "Will Maberry"
$ws.Range("C27").Value = "Will Maberry"

$ws.Range("A28").Value = "Index.html"
$ws.Range("B28").Value = "David Oyekola"
$ws.Range("C28").Value = "David Oyekola"

$ws.Range("A29").Value = "WordBankGame.java"
$ws.Range("B29").Value = "David Oyekola"
$ws.Range("C29").Value = "Didn't end up using"

$ws.Range("A30").Value = "Deployed to cse3310.org"
$ws.Range("B30").Value = "David Oyekola"
$ws.Range("C30").Value = "David Oyekola"

$ws.Range("A31").Value = "App.java (and associated unit tests)"
$ws.Range("B31").Value = "Grace Daily"
$ws.Range("C31").Value = "Grace Daily and Will Maberry"

$ws.Range("A32").Value = "ServerEvent.java (and associated unit tests)"
$ws.Range("B32").Value = "Grace Daily"
$ws.Range("C32").Value = "Grace Daily"

$ws.Range("A33").Value = "UserAuthentication.java (and associated unit tests)"
$ws.Range("B33").Value = "Ammar Rafiq"
$ws.Range("C33").Value = "Ammar Rafiq (no unit tests)"

$ws.Range("A34").Value = "UserEvent.java (and associated unit tests)"
$ws.Range("B34").Value = "Ammar Rafiq"
$ws.Range("C34").Value = "Will Maberry"

$ws.Range("A35").Value = "Lobby.java (and associated unit tests)"
$ws.Range("B35").Value = "Subodh Neupane"
$ws.Range("C35").Value = "Subodh Neupane (no unit tests)"

$ws.Range("A36").Value = "Leaderboard.java"
$ws.Range("B36").Value = "Subodh Neupane"
$ws.Range("C36").Value = "Subodh Neupane (no unit tests)"

$ws.Range("A37").Value = "Players.java"
$ws.Range("B37").Value = "Cody Mercer"
$ws.Range("C37").Value = "No work done and no communication"

$ws.Range("A38").Value = "Statistics.java"
$ws.Range("B38").Value = "Cody Mercer"
$ws.Range("C38").Value = "No work done and no communication"

# ---------------------------------------------------------------------
# 3) Restore the selection to match the edited workbook (C19:C29,
#    active cell C19).
# ---------------------------------------------------------------------
$ws.Range("C19:C29").Select() | Out-Null
